# Updated symbol list on Mon Dec 26 10:49:26 UTC 2022 with GitHub Actions
#
# The crypto-ranking table (rows 4-23) shifted down by one: "LEO" (which
# used to sit at row 23) is now listed right after "FTXToken" at row 4,
# and "One" (which used to sit further down) moved up to row 11; every
# other coin between them shifts down by one row accordingly. Prices
# (column D) and the "Volume(1h)" rank strings (column E) were refreshed
# to match. A handful of other price cells (row 2, 24, 26, 40-50) were
# also refreshed with new quotes, independent of the re-ordering.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Text-valued cells (coin name / link / volume-rank string): plain
#     assignment is safe because Excel won't try to reinterpret these as
#     numbers. ---
$textCells = @{
    "B4"  = "LEO"
    "C4"  = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
    "E4"  = "3LEOLEO"

    "B5"  = "HuobiToken"
    "C5"  = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
    "E5"  = "4HuobiTokenHT"

    "B6"  = "Cronos"
    "C6"  = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
    "E6"  = "5CronosCRO"

    "B7"  = "GateToken"
    "C7"  = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
    "E7"  = "6GateTokenGT"

    "B8"  = "KuCoinToken"
    "C8"  = "https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs"
    "E8"  = "7KuCoinTokenKCS"

    "B9"  = "MXToken"
    "C9"  = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
    "E9"  = "8MXTokenMX"

    "B10" = "FTXToken"
    "C10" = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
    "E10" = "9FTXTokenFTT"

    "B11" = "One"
    "C11" = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
    "E11" = "10OneONEBestin24h"

    "B12" = "WazirX"
    "C12" = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
    "E12" = "11WazirXWRX"

    "B13" = "MandalaExchangeToken"
    "C13" = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
    "E13" = "12MandalaExchangeTokenMDX"

    "B14" = "LiechtensteinCryptoassetsExchange"
    "C14" = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
    "E14" = "13LiechtensteinCryptoassetsExchangeLCXWorstin24h"

    "B15" = "BitrueCoin"
    "C15" = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
    "E15" = "14BitrueCoinBTR"

    "B16" = "BitMartToken"
    "C16" = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
    "E16" = "15BitMartTokenBMX"

    "B17" = "MCDex"
    "C17" = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
    "E17" = "16MCDexMCB"

    "B18" = "BitForexToken"
    "C18" = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
    "E18" = "17BitForexTokenBF"

    "B19" = "CoinExToken"
    "C19" = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
    "E19" = "18CoinExTokenCET"

    "B20" = "TigerCash"
    "C20" = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
    "E20" = "19TigerCashTCH"

    "B21" = "HotbitToken"
    "C21" = "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
    "E21" = "20HotbitTokenHTB"

    "B22" = "BitKan"
    "C22" = "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
    "E22" = "21BitKanKAN"

    "B23" = "NitroEx"
    "C23" = "https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx"
    "E23" = "22NitroExNTX"
}

foreach ($addr in $textCells.Keys) {
    $ws.Range($addr).Value = $textCells[$addr]
}

# --- Price cells (column D): these are numeric-looking strings stored as
#     *text* (to preserve exact formatting such as trailing zeros, e.g.
#     "2.150" rather than 2.15). A plain `.Value = "2.150"` assignment
#     gets auto-coerced by Excel into a binary double (losing the
#     trailing zero, and also introducing floating-point noise such as
#     243.38999999999999). Forcing the cell to Text format first makes
#     Excel store the literal characters instead; clearing the format
#     afterwards drops the leftover "number stored as text" styling so
#     the cell ends up with no explicit style, same as its neighbours. ---
$priceCells = @{
    "D2"  = "243.39"

    "D4"  = "3.610"
    "D5"  = "5.393"
    "D6"  = "0.05918"
    "D7"  = "3.453"
    "D8"  = "6.545"
    "D9"  = "0.8114"
    "D10" = "0.9106"
    "D11" = "0.01122"
    "D12" = "0.1412"
    "D13" = "0.07384"
    "D14" = "0.03275"
    "D15" = "0.03070"
    "D16" = "0.09352"
    "D17" = "3.848"
    "D18" = "0.001561"
    "D19" = "0.04679"
    "D20" = "0.006114"
    "D21" = "0.004990"
    "D22" = "0.0009822"
    "D23" = "0.00008603"

    "D24" = "2.150"
    "D26" = "0.1323"

    "D40" = "0.03958"
    "D41" = "0.006200"
    "D42" = "0.1076"
    "D43" = "0.003001"
    "D44" = "0.008147"
    "D45" = "0.00005156"
    "D47" = "0.8993"
    "D48" = "0.002338"
    "D49" = "0.00002101"
    "D50" = "0.0002001"
}

foreach ($addr in $priceCells.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $priceCells[$addr]
    $cell.ClearFormats()
}
